$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.030.89'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '2.956.28'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '379.28'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.46'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('E7').Value = '  +1.96%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.585'
$ws.Range('E9').Value = '  +0.80%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('E12').Value = '  +2.22%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '3.429.43'
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.42'
$ws.Range('E14').Value = '  +2.95%  '
$ws.Range('B15').Value = 'Uniswap'
$ws.Range('C15').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '12.41'
$ws.Range('E15').Value = '  +74.67%  '
$ws.Range('E16').Value = '  +5.85%  '
$ws.Range('D17').Value = '2.951.57'
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('E18').Value = '  +4.54%  '
$ws.Range('D19').Value = '51.076.57'
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('E20').Value = '  -2.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.40'
$ws.Range('E21').Value = '  -0.32%  '
$ws.Range('E22').Value = '  +0.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.35'
$ws.Range('E23').Value = '  +17.29%  '
$ws.Range('E24').Value = '  +2.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '267.35'
$ws.Range('E25').Value = '  +2.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.97'
$ws.Range('E26').Value = '  -2.54%  '
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('E28').Value = '  -0.82%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '25.78'
$ws.Range('E29').Value = '  +1.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.99'
$ws.Range('E30').Value = '  -6.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.108'
$ws.Range('E31').Value = '  -3.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.53'
$ws.Range('E32').Value = '  +7.88%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '50.79'
$ws.Range('E33').Value = '  +0.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '34.03'
$ws.Range('E34').Value = '  +0.88%  '
$ws.Range('E35').Value = '  +2.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0435'
$ws.Range('E36').Value = '  -3.58%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.22'
$ws.Range('E38').Value = '  +8.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.74'
$ws.Range('E39').Value = '  +3.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.116'
$ws.Range('E40').Value = '  +2.02%  '
$ws.Range('E41').Value = '  +3.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.49'
$ws.Range('E42').Value = '  -2.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '118.60'
$ws.Range('E43').Value = '  -1.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.57'
$ws.Range('E44').Value = '  +11.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.58'
$ws.Range('E45').Value = '  +2.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.03'
$ws.Range('E46').Value = '  -0.85%  '
$ws.Range('D47').Value = '2.040.52'
$ws.Range('E47').Value = '  +2.04%  '
$ws.Range('E48').Value = '  -1.01%  '
$ws.Range('E49').Value = '  -4.70%  '
$ws.Range('E50').Value = '  -6.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.34'
$ws.Range('E51').Value = '  +7.09%  '
